# added extent report in context with threading
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# Insert a new "arnav" row right after the existing "jyoti" row (row 4)
# in the AddCustomerTest block, pushing the OpenAccountTest block down.
$ws.Rows.Item(5).EntireRow.Insert()

# Append a new "arnav k" row at the bottom of the OpenAccountTest block first,
# so that the shared-string table records "arnav k" before "arnav".
$ws.Cells.Item(11, 1).Value = "Y"
$ws.Cells.Item(11, 2).Value = "arnav k"
$ws.Cells.Item(11, 3).Value = "Dollar"
$ws.Cells.Item(11, 4).Value = "firefox"

$ws.Cells.Item(5, 1).Value = "Y"
$ws.Cells.Item(5, 2).Value = "arnav"
$ws.Cells.Item(5, 3).Value = "k"
$ws.Cells.Item(5, 4).Value = "X7878"
$ws.Cells.Item(5, 5).Value = "firefox"

$ws.Range("C8").Select()
